# Cloudflare DNS export auto-update.
# A new DNS record (api.gpfree.org) was added by Cloudflare between the
# "ai.irrazionale.org" (row 122) and "freeundergroundtekno.org" (row 121)
# entries, pushing every following row down by one (old row 211 -> new row
# 212) and growing the sheet's used range from A1:P211 to A1:P212.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at 122 - this shifts rows 122..211 down to 123..212
# (and bumps the sheet dimension to A1:P212) exactly like Excel's own
# "Insert Sheet Rows" command would.
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row with the new DNS record's data.
$ws.Cells.Item(122, 1).Value = "edbbb21aeca4bf8ac037fabbcbeb403a"
$ws.Cells.Item(122, 2).Value = "api.gpfree.org"
$ws.Cells.Item(122, 3).Value = "A"
$ws.Cells.Item(122, 4).Value = "100.98.112.23"
$ws.Cells.Item(122, 5).Value = $false
$ws.Cells.Item(122, 6).Value = $false
$ws.Cells.Item(122, 7).Value = 1
$ws.Cells.Item(122, 8).Value = "{}"
$ws.Cells.Item(122, 9).Value = "{}"
$ws.Cells.Item(122, 11).Value = "[]"
$ws.Cells.Item(122, 12).Value = "2026-02-20T08:59:49.153859Z"
$ws.Cells.Item(122, 13).Value = "2026-02-20T08:59:49.153859Z"
